$wb = $excel.ActiveWorkbook

# --- "Port" sheet: update row 22, add row 23 ---
$port = $wb.Worksheets.Item("Port")

$port.Range("D22").Value = 'Issue with compiling common header files'
$port.Range("E22").Value = 'Don’t compile common header files as C++ (/TP)'

$port.Range("A23").Value = 'step7'
$port.Range("B23").Value = 'Compile DS (Diagnostic Server code)'
$port.Range("C23").Value = 'Error free for DS server'
$port.Range("D23").Value = 'LF_DIAGNOSTICS_PORT'
$port.Range("E21").Copy()
$port.Range("D23").PasteSpecial(-4122)
$port.Range("E23").Value = 'This is cause common code calls functions like DS_LOG_ERROR_2, which in coreclr has the constant defined'
$port.Range("F23").Value = 'Tried getting to these constants via stressLog.h but trying to punt for now'

$port.Columns.Item(5).AutoFit()

# --- "Port_Issues" sheet: add rows 26-37 ---
$issues = $wb.Worksheets.Item("Port_Issues")

$issues.Range("E26").Value = '''''-IC:\Work\Core\CurrentWork\runtime\artifacts\obj\coreclr\windows.x64.Release\nativeaot\Runtime\Full '
$issues.Range("E26").Style = "Normal"
$issues.Range("E27").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\Full '
$issues.Range("E27").Style = "Normal"
$issues.Range("E28").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\native '
$issues.Range("E28").Style = "Normal"
$issues.Range("E29").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\pal\prebuilt\inc '
$issues.Range("E29").Style = "Normal"
$issues.Range("E30").Value = '''''-IC:\Work\Core\CurrentWork\runtime\artifacts\obj '
$issues.Range("E30").Style = "Normal"
$issues.Range("E31").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\native\eventpipe '
$issues.Range("E31").Style = "Normal"
$issues.Range("E32").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\inc '
$issues.Range("E32").Style = "Normal"
$issues.Range("E33").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\. '
$issues.Range("E33").Style = "Normal"
$issues.Range("E34").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\..\..\gc '
$issues.Range("E34").Style = "Normal"
$issues.Range("E35").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\..\..\gc\env '
$issues.Range("E35").Style = "Normal"
$issues.Range("E36").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\windows '
$issues.Range("E36").Style = "Normal"
$issues.Range("E37").Value = '''''-IC:\Work\Core\CurrentWork\runtime\src\coreclr\nativeaot\Runtime\amd64 '
$issues.Range("E37").Style = "Normal"

$issues.Columns.Item(5).AutoFit()

# --- selections / active sheet bookkeeping ---
$tables = $wb.Worksheets.Item("Port-Tables")
$tables.Range("B30").Select()

$issues.Range("E32").Select()

$port.Range("D23").Select()

